$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 25/26: Kaspa and Fetch.AI swap positions (with new Volume values)
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D25") "1.63"
$ws.Range("E25").Value = "  +9.12%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D26") "0.170"
$ws.Range("E26").Value = "  -2.07%  "

Set-TextValue $ws.Range("D2") "63.120.89"
$ws.Range("E2").Value = "  +0.25%  "

Set-TextValue $ws.Range("D3") "2.543.82"
$ws.Range("E3").Value = "  +3.17%  "

Set-TextValue $ws.Range("D5") "568.58"
$ws.Range("E5").Value = "  +0.66%  "

Set-TextValue $ws.Range("D6") "147.38"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("E7").Value = "  +0.03%  "

Set-TextValue $ws.Range("D8") "0.586"
$ws.Range("E8").Value = "  -0.40%  "

Set-TextValue $ws.Range("D9") "2.541.68"
$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("E13").Value = "  +0.27%  "

Set-TextValue $ws.Range("D14") "27.50"
$ws.Range("E14").Value = "  +4.51%  "

Set-TextValue $ws.Range("D15") "3.001.46"
$ws.Range("E15").Value = "  +3.28%  "

Set-TextValue $ws.Range("D16") "63.081.70"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("E17").Value = "  +2.00%  "

Set-TextValue $ws.Range("D18") "2.545.26"
$ws.Range("E18").Value = "  +3.17%  "

Set-TextValue $ws.Range("D19") "11.48"
$ws.Range("E19").Value = "  +2.12%  "

Set-TextValue $ws.Range("D20") "336.13"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("E21").Value = "  +0.89%  "

Set-TextValue $ws.Range("D22") "6.76"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("E23").Value = "  +0.12%  "

Set-TextValue $ws.Range("D24") "65.30"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("E27").Value = "  +11.54%  "

Set-TextValue $ws.Range("D28") "8.46"
$ws.Range("E28").Value = "  +5.06%  "

Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.09%  "

Set-TextValue $ws.Range("D30") "7.34"
$ws.Range("E30").Value = "  +7.62%  "

$ws.Range("E31").Value = "  +2.88%  "

Set-TextValue $ws.Range("D32") "1.85"
$ws.Range("E32").Value = "  +0.68%  "

Set-TextValue $ws.Range("D33") "175.94"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("E34").Value = "  +3.96%  "

Set-TextValue $ws.Range("D35") "411.23"
$ws.Range("E35").Value = "  +12.51%  "

Set-TextValue $ws.Range("D36") "0.398"
$ws.Range("E36").Value = "  +0.41%  "

Set-TextValue $ws.Range("D37") "19.02"
$ws.Range("E37").Value = "  +0.92%  "

Set-TextValue $ws.Range("D38") "4.39"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("E41").Value = "  -0.02%  "

Set-TextValue $ws.Range("D42") "39.29"
$ws.Range("E42").Value = "  -2.92%  "

Set-TextValue $ws.Range("D43") "153.38"
$ws.Range("E43").Value = "  +2.35%  "

Set-TextValue $ws.Range("D44") "3.77"
$ws.Range("E44").Value = "  +2.07%  "

Set-TextValue $ws.Range("D45") "20.94"
$ws.Range("E45").Value = "  +2.01%  "

Set-TextValue $ws.Range("D46") "0.604"
$ws.Range("E46").Value = "  +1.09%  "

Set-TextValue $ws.Range("D47") "0.0962"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("E48").Value = "  +6.01%  "

Set-TextValue $ws.Range("D49") "0.0521"
$ws.Range("E49").Value = "  +1.21%  "

Set-TextValue $ws.Range("D50") "18.30"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("E51").Value = "  +2.36%  "
